$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three data rows (2, 3, 4) have their Fecha / Volumen / Precio minimo /
# Precio maximo / Precio promedio ponderado / Origen / Precio $/Kg values
# cyclically rotated: row2 <- row3, row3 <- row4, row4 <- row2.
# Capture the "before" values first, then write them back in rotated order.

$cols = @("D", "J", "K", "L", "M", "O", "P")

$row2 = @{}
$row3 = @{}
$row4 = @{}
foreach ($col in $cols) {
    $row2[$col] = $ws.Range("${col}2").Value2
    $row3[$col] = $ws.Range("${col}3").Value2
    $row4[$col] = $ws.Range("${col}4").Value2
}

foreach ($col in $cols) {
    $ws.Range("${col}2").Value2 = $row3[$col]
    $ws.Range("${col}3").Value2 = $row4[$col]
    $ws.Range("${col}4").Value2 = $row2[$col]
}
